$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new bibliography entry ("Prikryl, P., Themens, D. R., Chum, J.,
#    Chakraborty, S., Gillies, R. G., and Weygand, J. M.: ... 2025.") right
#    before the existing "Chakraborty, S., Qian, L., ..." entry. Everything
#    that previously followed shifts down by one list item, matching the
#    target diff exactly.
# ---------------------------------------------------------------------------

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t.StartsWith("Chakraborty") -and $t -like "*G-condition*" -and $t -like "*Great American Eclipse*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'Chakraborty, S., Qian, L...' bibliography paragraph"
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $targetPara.Range
$targetRange.InsertParagraphBefore()

# Re-fetch the freshly created (now empty) paragraph that precedes the
# original one (it now occupies the same index the original paragraph used
# to have), and fill it in with the exact run/formatting structure from the
# target XML (a plain run, a spell-check-flagged run for "Themens", a bold
# run for "Chakraborty, S.", and a trailing plain run).
$newPara = $d.Paragraphs.Item($targetIndex)
$newRange = $newPara.Range

$prikrylXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Prikryl, P., </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Themens</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, D. R., Chum, J., </w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Chakraborty, S.</w:t></w:r>' +
    '<w:r><w:t>, Gillies, R. G., and Weygand, J. M.: Observations of traveling ionospheric disturbances driven by gravity waves from sources in the upper and lower atmosphere, Ann. Geophys., 43, 511' +
    [char]0x2013 +
    '534, https://doi.org/10.5194/angeo-43-511-2025, 2025.</w:t></w:r>' +
    '</w:p>'

$newRange.InsertXML($prikrylXml)

# ---------------------------------------------------------------------------
# 2) Merge the three runs ("...magnetosphere, Ann. " / "Geophys" / "., 40,
#    619-639, ") describing the *other*, pre-existing Prikryl et al. (2022)
#    reference into a single run, dropping the spell-check markers around
#    "Geophys" in the process (the visible text is unchanged).
# ---------------------------------------------------------------------------

$find = $d.Content
$found = $find.Find.Execute("Ann. Geophys., 40, 619-639, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Ann. Geophys., 40, 619-639, ' run sequence"
}
$mergeStart = $find.Start
$mergeEnd = $find.End

# Setting identical text is a no-op for the underlying run structure, so
# stage a distinct placeholder first and then restore the real text; that
# forces the three runs (and the "Geophys" proofErr wrapper) to collapse
# into a single plain run, matching the target markup exactly.
$mergeRange = $d.Range($mergeStart, $mergeEnd)
$mergeRange.Text = "TEMP_PLACEHOLDER_FOR_MERGE"
$mergeRange2 = $d.Range($mergeStart, $mergeStart + "TEMP_PLACEHOLDER_FOR_MERGE".Length)
$mergeRange2.Text = "Ann. Geophys., 40, 619-639, "

Write-Host "Edit complete."
